$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug List")

# New bug entries added to the "Bug List" sheet, all found in version 1.0.6
$ws.Range("A13").Value = "Dummy account appears in the order status widget when creating a new order"
$ws.Range("B13").Value = "1.0.6"

$ws.Range("A14").Value = "when customer is on hold, prevent orders being created"
$ws.Range("B14").Value = "1.0.6"

$ws.Range("A15").Value = "remove the bin description from the loader sheet"
$ws.Range("B15").Value = "1.0.6"

$ws.Range("A16").Value = "create a separate sheet for additivies"
$ws.Range("B16").Value = "1.0.6"

$ws.Range("A17").Value = "change the deliery page to selec the truck first"
$ws.Range("B17").Value = "1.0.6"

# Move the active selection/tab to the Bug List sheet, cell C17
$ws.Range("C17").Select() | Out-Null
